{"js": "// Update \"Traits table (new citation format).docx\":\n// A new citation (55) was inserted into several citation-number cells, the\n// phrase \"carbon isotope composition\" was reworded to \"carbon isotope ratio\",\n// and one citation list was reordered so \"51\" comes first.\n//\n// Each target cell is located by its exact, unique paragraph text (the\n// citation-number lists in this table are unique per row), then the text is\n// replaced in place so the surrounding table cell/paragraph formatting is\n// preserved.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Map of the paragraph's current exact text -> its new exact text.\nconst replacements = new Map([\n  [\"1, 7, 2, 3, 4, 6\", \"1, 55, 7, 2, 3, 4, 6\"],\n  [\"7, 29, 30, 32, 31, 9\", \"55, 7, 29, 30, 32, 31, 9\"],\n  [\"15, 7, 29, 30, 32, 34\", \"55, 15, 7, 29, 30, 32, 34\"],\n  [\"15, 36, 1, 37   \", \"55, 15, 36, 1, 37   \"],\n  [\"15, 35, 1\", \"55, 15, 35, 1\"],\n  [\"carbon isotope composition\", \"carbon isotope ratio\"],\n  [\"7, 43, 31\", \"55, 7, 43, 31\"],\n  [\"30, 44, 51\", \"51, 30, 44\"],\n]);\n\nlet remaining = replacements.size;\nfor (let i = 0; i < paragraphs.items.length && remaining > 0; i++) {\n  const para = paragraphs.items[i];\n  const newText = replacements.get(para.text);\n  if (newText !== undefined) {\n    para.insertText(newText, Word.InsertLocation.replace);\n    replacements.delete(para.text);\n    remaining--;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update \"Traits table (new citation format).docx\":\n# A new citation (55) was inserted into several citation-number cells, the\n# phrase \"carbon isotope composition\" was reworded to \"carbon isotope ratio\",\n# and one citation list was reordered so \"51\" comes first.\n#\n# Each target cell is located by its exact, unique paragraph text (the\n# citation-number lists in this table are unique per row), then the text is\n# replaced in place via a Find/Replace scoped to that paragraph's range so\n# the surrounding table cell/paragraph formatting is preserved.\n\n$d = $word.ActiveDocument\n\n$replacements = @{\n    \"1, 7, 2, 3, 4, 6\"            = \"1, 55, 7, 2, 3, 4, 6\"\n    \"7, 29, 30, 32, 31, 9\"        = \"55, 7, 29, 30, 32, 31, 9\"\n    \"15, 7, 29, 30, 32, 34\"       = \"55, 15, 7, 29, 30, 32, 34\"\n    \"15, 36, 1, 37   \"            = \"55, 15, 36, 1, 37   \"\n    \"15, 35, 1\"                   = \"55, 15, 35, 1\"\n    \"carbon isotope composition\"  = \"carbon isotope ratio\"\n    \"7, 43, 31\"                   = \"55, 7, 43, 31\"\n    \"30, 44, 51\"                  = \"51, 30, 44\"\n}\n\n$remaining = $replacements.Count\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count -and $remaining -gt 0; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $rng = $p.Range\n    # Drop the trailing paragraph mark / cell-end mark so the comparison is\n    # against the visible cell text only.\n    $rng.MoveEnd(1, -1) | Out-Null\n    $text = $rng.Text\n\n    if ($replacements.ContainsKey($text)) {\n        $newText = $replacements[$text]\n        $rng.Find.Execute($text, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n        $replacements.Remove($text)\n        $remaining--\n    }\n}\n"}
